# Updated: po 26. 04. 2021
# Revises the AgTests (F) / AgPosit (G) columns for rows 307-414 with
# corrected antigen-testing figures, and appends two new daily rows
# (415, 416) for 2021-04-23 and 2021-04-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing cells (columns F and G) for rows 307-414 ---
$ws.Range("F307").Value = 75444
$ws.Range("G307").Value = 6335
$ws.Range("F326").Value = 418593
$ws.Range("F334").Value = 192789
$ws.Range("G334").Value = 3500
$ws.Range("F335").Value = 150596
$ws.Range("G335").Value = 3781
$ws.Range("F338").Value = 221555
$ws.Range("F348").Value = 232777
$ws.Range("F353").Value = 723698
$ws.Range("F356").Value = 160044
$ws.Range("F358").Value = 158794
$ws.Range("F360").Value = 750085
$ws.Range("G360").Value = 5142
$ws.Range("F362").Value = 229210
$ws.Range("F363").Value = 188694
$ws.Range("F364").Value = 168396
$ws.Range("F365").Value = 184874
$ws.Range("F366").Value = 339399
$ws.Range("F367").Value = 766852
$ws.Range("G367").Value = 3920
$ws.Range("F368").Value = 344996
$ws.Range("G368").Value = 2289
$ws.Range("F369").Value = 234663
$ws.Range("G369").Value = 2599
$ws.Range("F370").Value = 180742
$ws.Range("F371").Value = 160129
$ws.Range("G371").Value = 1957
$ws.Range("F372").Value = 179033
$ws.Range("F374").Value = 773680
$ws.Range("G374").Value = 3420
$ws.Range("F375").Value = 350178
$ws.Range("G375").Value = 1852
$ws.Range("F376").Value = 221387
$ws.Range("G376").Value = 2224
$ws.Range("F377").Value = 176981
$ws.Range("G377").Value = 1826
$ws.Range("F378").Value = 157258
$ws.Range("G378").Value = 1550
$ws.Range("F379").Value = 179971
$ws.Range("G379").Value = 1619
$ws.Range("F380").Value = 344836
$ws.Range("G380").Value = 2024
$ws.Range("F381").Value = 746054
$ws.Range("G381").Value = 2695
$ws.Range("F382").Value = 356983
$ws.Range("F383").Value = 220777
$ws.Range("G383").Value = 1766
$ws.Range("F384").Value = 172043
$ws.Range("G384").Value = 1513
$ws.Range("F385").Value = 150902
$ws.Range("F386").Value = 182753
$ws.Range("G386").Value = 1361
$ws.Range("F387").Value = 351511
$ws.Range("G387").Value = 1666
$ws.Range("F388").Value = 729618
$ws.Range("G388").Value = 2203
$ws.Range("F389").Value = 353665
$ws.Range("G389").Value = 1305
$ws.Range("F390").Value = 219769
$ws.Range("G390").Value = 1474
$ws.Range("F391").Value = 177297
$ws.Range("G391").Value = 1208
$ws.Range("F392").Value = 220981
$ws.Range("G392").Value = 1217
$ws.Range("F393").Value = 307442
$ws.Range("G393").Value = 1232
$ws.Range("F394").Value = 166346
$ws.Range("G394").Value = 633
$ws.Range("F395").Value = 750877
$ws.Range("G395").Value = 1956
$ws.Range("F397").Value = 108084
$ws.Range("G397").Value = 640
$ws.Range("F398").Value = 298560
$ws.Range("F399").Value = 201110
$ws.Range("G399").Value = 970
$ws.Range("F400").Value = 150311
$ws.Range("G400").Value = 759
$ws.Range("F401").Value = 273338
$ws.Range("G401").Value = 934
$ws.Range("F402").Value = 716838
$ws.Range("G402").Value = 1385
$ws.Range("F403").Value = 351375
$ws.Range("F404").Value = 224415
$ws.Range("G404").Value = 906
$ws.Range("F405").Value = 173948
$ws.Range("G405").Value = 695
$ws.Range("F406").Value = 170765
$ws.Range("G406").Value = 679
$ws.Range("F407").Value = 157304
$ws.Range("F408").Value = 303125
$ws.Range("G408").Value = 839
$ws.Range("F409").Value = 698729
$ws.Range("G409").Value = 1001
$ws.Range("F410").Value = 351563
$ws.Range("G410").Value = 620
$ws.Range("F411").Value = 225120
$ws.Range("G411").Value = 822
$ws.Range("F412").Value = 175664
$ws.Range("G412").Value = 642
$ws.Range("F413").Value = 148664
$ws.Range("G413").Value = 658
$ws.Range("F414").Value = 144902
$ws.Range("G414").Value = 550

# --- Append new rows 415 and 416 ---
$ws.Range("A415").Value = 44309
$ws.Range("B415").Value = 379476
$ws.Range("C415").Value = 14312
$ws.Range("D415").Value = 841
$ws.Range("E415").Value = 11458
$ws.Range("F415").Value = 279657
$ws.Range("G415").Value = 757
$ws.Range("A416").Value = 44310
$ws.Range("B416").Value = 379911
$ws.Range("C416").Value = 30594
$ws.Range("D416").Value = 435
$ws.Range("E416").Value = 11495
$ws.Range("F416").Value = 554958
$ws.Range("G416").Value = 1169
